$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.362.25"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.873.12"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.015"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +1.36%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "238.43"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.015"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4704"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2831"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.88%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06414"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "18.18"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.888.35"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07591"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.82%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "94.87"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +11.49%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.985"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6431"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "288.82"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +18.58%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.484.25"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.008"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.187.20"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.000007379"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.39%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.012"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.037"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "165.28"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.123"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.73%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.28"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.88%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.923"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.32%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.1079"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.91%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.350"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.49%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.043"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.782"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.20%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04931"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.59%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7247"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.40%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.114"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.52%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.749"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.28%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01942"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.697"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.993"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.49%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.8645"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "106.69"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.34%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.012"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.11%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.606"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.4101"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "65.03"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.02%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "7.057"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.13%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.984"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.06%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1190"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.44%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "34.10"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05620"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3752"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
